$d = $word.ActiveDocument

# Delete the paragraphs that are no longer needed now that the
# "add user" dialog exists:
#   - "Cloud funtions: crear usuario en base de datos" (paragraph 4)
#   - "Selector del menú al volver de settings" (paragraph 2)
# Deleting from the highest paragraph index down avoids index shifting.
$d.Paragraphs.Item(4).Range.Delete()
$d.Paragraphs.Item(2).Range.Delete()
